$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 117.6  # H6: 131.11111 -> 117.6
$ws.Cells.Item(6, 9).Value = 75.22221999999999  # I6: 85.125 -> 75.22221999999999
$ws.Cells.Item(6, 11).Value = 225.66666  # K6: 255.375 -> 225.66666
$ws.Cells.Item(6, 13).Value = -113.66666  # M6: -143.375 -> -113.66666
$ws.Cells.Item(98, 8).Value = 1528.5333  # H98: 1578.4828 -> 1528.5333
$ws.Cells.Item(98, 9).Value = 1562.3462  # I98: 1621.64 -> 1562.3462
$ws.Cells.Item(98, 11).Value = 1562.3462  # K98: 1621.64 -> 1562.3462
$ws.Cells.Item(98, 13).Value = -64.34619999999995  # M98: -123.6400000000001 -> -64.34619999999995
$ws.Cells.Item(103, 8).Value = 1550.5  # H103: 1874.75 -> 1550.5
$ws.Cells.Item(103, 9).Value = 1325.75  # I103: 1749.5 -> 1325.75
$ws.Cells.Item(103, 11).Value = 3977.25  # K103: 5248.5 -> 3977.25
$ws.Cells.Item(103, 13).Value = -3391.25  # M103: -4662.5 -> -3391.25
$ws.Cells.Item(122, 8).Value = 1528.5333  # H122: 1578.4828 -> 1528.5333
$ws.Cells.Item(122, 9).Value = 1562.3462  # I122: 1621.64 -> 1562.3462
$ws.Cells.Item(122, 11).Value = 4687.0386  # K122: 4864.92 -> 4687.0386
$ws.Cells.Item(122, 13).Value = -2237.0386  # M122: -2414.92 -> -2237.0386
$ws.Cells.Item(125, 8).Value = 33335064  # H125: 41668580 -> 33335064
$ws.Cells.Item(125, 9).Value = 1647.25  # I125: 2033.3334 -> 1647.25
$ws.Cells.Item(125, 10).Value = 45456304  # J125: 55557428 -> 45456304
$ws.Cells.Item(125, 11).Value = 14825.25  # K125: 18300.0006 -> 14825.25
$ws.Cells.Item(125, 12).Value = 409106736  # L125: 500016852 -> 409106736
$ws.Cells.Item(125, 13).Value = -12365.25  # M125: -15840.0006 -> -12365.25
$ws.Cells.Item(125, 14).Value = -409111656  # N125: -500021772 -> -409111656
$ws.Cells.Item(137, 8).Value = 2111.8057  # H137: 2079.027 -> 2111.8057
$ws.Cells.Item(137, 9).Value = 1328.9231  # I137: 1298.2142 -> 1328.9231
$ws.Cells.Item(137, 11).Value = 3986.7693  # K137: 3894.6426 -> 3986.7693
$ws.Cells.Item(137, 13).Value = -1436.7693  # M137: -1344.6426 -> -1436.7693

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 4266.6665  # H3: 4000.75 -> 4266.6665
$ws.Cells.Item(3, 9).Value = 0  # I3: 4000.75 -> 0
$ws.Cells.Item(3, 10).Value = 4266.6665  # J3: 0 -> 4266.6665
$ws.Cells.Item(3, 11).Value = 0  # K3: 4000.75 -> 0
$ws.Cells.Item(3, 12).Value = 4266.6665  # L3: 0 -> 4266.6665
$ws.Cells.Item(3, 13).ClearContents()  # M3: was -3885.75
$ws.Cells.Item(3, 14).Value = -4496.6665  # N3: None -> -4496.6665
$ws.Cells.Item(32, 8).Value = 4357  # H32: 5429.48 -> 4357
$ws.Cells.Item(32, 9).Value = 2095.4443  # I32: 3151.6328 -> 2095.4443
$ws.Cells.Item(32, 10).Value = 13998.368  # J32: 13998.523 -> 13998.368
$ws.Cells.Item(32, 11).Value = 2095.4443  # K32: 3151.6328 -> 2095.4443
$ws.Cells.Item(32, 12).Value = 13998.368  # L32: 13998.523 -> 13998.368
$ws.Cells.Item(32, 13).Value = -1808.4443  # M32: -2864.6328 -> -1808.4443
$ws.Cells.Item(32, 14).Value = -14572.368  # N32: -14572.523 -> -14572.368
$ws.Cells.Item(45, 8).Value = 4021.7856  # H45: 6136.5835 -> 4021.7856
$ws.Cells.Item(45, 9).Value = 2946.5386  # I45: 5058.091 -> 2946.5386
$ws.Cells.Item(45, 11).Value = 2946.5386  # K45: 5058.091 -> 2946.5386
$ws.Cells.Item(45, 13).Value = -2569.5386  # M45: -4681.091 -> -2569.5386
$ws.Cells.Item(61, 8).Value = 11467.344  # H61: 12414.793 -> 11467.344
$ws.Cells.Item(61, 9).Value = 16560.818  # I61: 18019.3 -> 16560.818
$ws.Cells.Item(61, 10).Value = 8799.333000000001  # J61: 9465.053 -> 8799.333000000001
$ws.Cells.Item(61, 11).Value = 16560.818  # K61: 18019.3 -> 16560.818
$ws.Cells.Item(61, 12).Value = 8799.333000000001  # L61: 9465.053 -> 8799.333000000001
$ws.Cells.Item(61, 13).Value = -16348.818  # M61: -17807.3 -> -16348.818
$ws.Cells.Item(61, 14).Value = -9223.333000000001  # N61: -9889.053 -> -9223.333000000001
$ws.Cells.Item(63, 8).Value = 7498.25  # H63: 9875 -> 7498.25
$ws.Cells.Item(63, 9).Value = 1995.8  # I63: 2000 -> 1995.8
$ws.Cells.Item(63, 10).Value = 11428.571  # J63: 12500 -> 11428.571
$ws.Cells.Item(63, 11).Value = 1995.8  # K63: 2000 -> 1995.8
$ws.Cells.Item(63, 12).Value = 11428.571  # L63: 12500 -> 11428.571
$ws.Cells.Item(63, 13).Value = -1309.8  # M63: -1314 -> -1309.8
$ws.Cells.Item(63, 14).Value = -12800.571  # N63: -13872 -> -12800.571
$ws.Cells.Item(66, 8).Value = 7498.25  # H66: 9875 -> 7498.25
$ws.Cells.Item(66, 9).Value = 1995.8  # I66: 2000 -> 1995.8
$ws.Cells.Item(66, 10).Value = 11428.571  # J66: 12500 -> 11428.571
$ws.Cells.Item(66, 11).Value = 9979  # K66: 10000 -> 9979
$ws.Cells.Item(66, 12).Value = 57142.855  # L66: 62500 -> 57142.855
$ws.Cells.Item(66, 13).Value = -6547  # M66: -6568 -> -6547
$ws.Cells.Item(66, 14).Value = -64006.855  # N66: -69364 -> -64006.855
$ws.Cells.Item(74, 8).Value = 1647257.5  # H74: 1691814 -> 1647257.5
$ws.Cells.Item(74, 9).Value = 2085413.5  # I74: 2157370.2 -> 2085413.5
$ws.Cells.Item(74, 11).Value = 2085413.5  # K74: 2157370.2 -> 2085413.5
$ws.Cells.Item(74, 13).Value = -2084539.5  # M74: -2156496.2 -> -2084539.5
$ws.Cells.Item(77, 8).Value = 1647257.5  # H77: 1691814 -> 1647257.5
$ws.Cells.Item(77, 9).Value = 2085413.5  # I77: 2157370.2 -> 2085413.5
$ws.Cells.Item(77, 11).Value = 10427067.5  # K77: 10786851 -> 10427067.5
$ws.Cells.Item(77, 13).Value = -10422699.5  # M77: -10782483 -> -10422699.5
$ws.Cells.Item(97, 8).Value = 641.9375  # H97: 703.6429000000001 -> 641.9375
$ws.Cells.Item(97, 9).Value = 573  # I97: 639 -> 573
$ws.Cells.Item(97, 11).Value = 573  # K97: 639 -> 573
$ws.Cells.Item(97, 13).Value = -77  # M97: -143 -> -77
$ws.Cells.Item(102, 8).Value = 5470.75  # H102: 2102.861 -> 5470.75
$ws.Cells.Item(102, 9).Value = 4828.2856  # I102: 2029.5161 -> 4828.2856
$ws.Cells.Item(102, 10).Value = 9968  # J102: 2557.6 -> 9968
$ws.Cells.Item(102, 11).Value = 4828.2856  # K102: 2029.5161 -> 4828.2856
$ws.Cells.Item(102, 12).Value = 9968  # L102: 2557.6 -> 9968
$ws.Cells.Item(102, 13).Value = -3206.2856  # M102: -407.5161000000001 -> -3206.2856
$ws.Cells.Item(102, 14).Value = -13212  # N102: -5801.6 -> -13212
$ws.Cells.Item(132, 8).Value = 442556.5  # H132: 442535.78 -> 442556.5
$ws.Cells.Item(132, 9).Value = 630979.6  # I132: 614396 -> 630979.6
$ws.Cells.Item(132, 10).Value = 6828  # J132: 7156.533 -> 6828
$ws.Cells.Item(132, 11).Value = 1892938.8  # K132: 1843188 -> 1892938.8
$ws.Cells.Item(132, 12).Value = 20484  # L132: 21469.599 -> 20484
$ws.Cells.Item(132, 13).Value = -1890408.8  # M132: -1840658 -> -1890408.8
$ws.Cells.Item(132, 14).Value = -25544  # N132: -26529.599 -> -25544
$ws.Cells.Item(136, 8).Value = 11467.344  # H136: 12414.793 -> 11467.344
$ws.Cells.Item(136, 9).Value = 16560.818  # I136: 18019.3 -> 16560.818
$ws.Cells.Item(136, 10).Value = 8799.333000000001  # J136: 9465.053 -> 8799.333000000001
$ws.Cells.Item(136, 11).Value = 49682.454  # K136: 54057.89999999999 -> 49682.454
$ws.Cells.Item(136, 12).Value = 26397.999  # L136: 28395.159 -> 26397.999
$ws.Cells.Item(136, 13).Value = -47132.454  # M136: -51507.89999999999 -> -47132.454
$ws.Cells.Item(136, 14).Value = -31497.999  # N136: -33495.159 -> -31497.999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2225.1853  # H20: 2334.68 -> 2225.1853
$ws.Cells.Item(20, 9).Value = 1910.579  # I20: 2100.1875 -> 1910.579
$ws.Cells.Item(20, 10).Value = 2972.375  # J20: 2751.5557 -> 2972.375
$ws.Cells.Item(20, 11).Value = 1910.579  # K20: 2100.1875 -> 1910.579
$ws.Cells.Item(20, 12).Value = 2972.375  # L20: 2751.5557 -> 2972.375
$ws.Cells.Item(20, 13).Value = -1663.579  # M20: -1853.1875 -> -1663.579
$ws.Cells.Item(20, 14).Value = -3466.375  # N20: -3245.5557 -> -3466.375
$ws.Cells.Item(22, 8).Value = 2008.1666  # H22: 2042.5714 -> 2008.1666
$ws.Cells.Item(22, 9).Value = 1609.8  # I22: 1716.3334 -> 1609.8
$ws.Cells.Item(22, 11).Value = 1609.8  # K22: 1716.3334 -> 1609.8
$ws.Cells.Item(22, 13).Value = -1436.8  # M22: -1543.3334 -> -1436.8
$ws.Cells.Item(86, 8).Value = 3368.2856  # H86: 2874.889 -> 3368.2856
$ws.Cells.Item(86, 9).Value = 1625  # I86: 1476.6666 -> 1625
$ws.Cells.Item(86, 10).Value = 4065.6  # J86: 5671.3335 -> 4065.6
$ws.Cells.Item(86, 11).Value = 1625  # K86: 1476.6666 -> 1625
$ws.Cells.Item(86, 12).Value = 4065.6  # L86: 5671.3335 -> 4065.6
$ws.Cells.Item(86, 13).Value = -502  # M86: -353.6666 -> -502
$ws.Cells.Item(86, 14).Value = -6311.6  # N86: -7917.3335 -> -6311.6
$ws.Cells.Item(89, 8).Value = 3368.2856  # H89: 2874.889 -> 3368.2856
$ws.Cells.Item(89, 9).Value = 1625  # I89: 1476.6666 -> 1625
$ws.Cells.Item(89, 10).Value = 4065.6  # J89: 5671.3335 -> 4065.6
$ws.Cells.Item(89, 11).Value = 8125  # K89: 7383.333000000001 -> 8125
$ws.Cells.Item(89, 12).Value = 20328  # L89: 28356.6675 -> 20328
$ws.Cells.Item(89, 13).Value = -2509  # M89: -1767.333000000001 -> -2509
$ws.Cells.Item(89, 14).Value = -31560  # N89: -39588.6675 -> -31560
$ws.Cells.Item(99, 8).Value = 2672.5454  # H99: 2895.8 -> 2672.5454
$ws.Cells.Item(99, 9).Value = 923.75  # I99: 992.8570999999999 -> 923.75
$ws.Cells.Item(99, 11).Value = 923.75  # K99: 992.8570999999999 -> 923.75
$ws.Cells.Item(99, 13).Value = 574.25  # M99: 505.1429000000001 -> 574.25
$ws.Cells.Item(105, 8).Value = 4398.5654  # H105: 4417.6816 -> 4398.5654
$ws.Cells.Item(105, 10).Value = 4612.933  # J105: 4658.2856 -> 4612.933
$ws.Cells.Item(105, 12).Value = 4612.933  # L105: 4658.2856 -> 4612.933
$ws.Cells.Item(105, 14).Value = -8106.933  # N105: -8152.2856 -> -8106.933
$ws.Cells.Item(134, 8).Value = 543146.2  # H134: 555237.25 -> 543146.2
$ws.Cells.Item(134, 9).Value = 656098.25  # I134: 673856.3 -> 656098.25
$ws.Cells.Item(134, 11).Value = 1968294.75  # K134: 2021568.9 -> 1968294.75
$ws.Cells.Item(134, 13).Value = -1965759.75  # M134: -2019033.9 -> -1965759.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 516815.22  # H58: 563631.4 -> 516815.22
$ws.Cells.Item(58, 9).Value = 687667.5  # I58: 687660.5600000001 -> 687667.5
$ws.Cells.Item(58, 10).Value = 4258.3335  # J58: 5500 -> 4258.3335
$ws.Cells.Item(58, 11).Value = 687667.5  # K58: 687660.5600000001 -> 687667.5
$ws.Cells.Item(58, 12).Value = 4258.3335  # L58: 5500 -> 4258.3335
$ws.Cells.Item(58, 13).Value = -687464.5  # M58: -687457.5600000001 -> -687464.5
$ws.Cells.Item(58, 14).Value = -4664.3335  # N58: -5906 -> -4664.3335
$ws.Cells.Item(62, 8).Value = 3583.111  # H62: 3821.2856 -> 3583.111
$ws.Cells.Item(62, 9).Value = 3199.6  # I62: 3499.6667 -> 3199.6
$ws.Cells.Item(62, 11).Value = 3199.6  # K62: 3499.6667 -> 3199.6
$ws.Cells.Item(62, 13).Value = -2575.6  # M62: -2875.6667 -> -2575.6
$ws.Cells.Item(65, 8).Value = 3583.111  # H65: 3821.2856 -> 3583.111
$ws.Cells.Item(65, 9).Value = 3199.6  # I65: 3499.6667 -> 3199.6
$ws.Cells.Item(65, 11).Value = 15998  # K65: 17498.3335 -> 15998
$ws.Cells.Item(65, 13).Value = -12878  # M65: -14378.3335 -> -12878
$ws.Cells.Item(107, 8).Value = 1084.7059  # H107: 1114.1111 -> 1084.7059
$ws.Cells.Item(107, 9).Value = 883.9231  # I107: 991.8182 -> 883.9231
$ws.Cells.Item(107, 10).Value = 1737.25  # J107: 1306.2858 -> 1737.25
$ws.Cells.Item(107, 11).Value = 883.9231  # K107: 991.8182 -> 883.9231
$ws.Cells.Item(107, 12).Value = 1737.25  # L107: 1306.2858 -> 1737.25
$ws.Cells.Item(107, 13).Value = 1036.0769  # M107: 928.1818 -> 1036.0769
$ws.Cells.Item(107, 14).Value = -5577.25  # N107: -5146.2858 -> -5577.25
$ws.Cells.Item(108, 8).Value = 66547.2  # H108: 70347.2 -> 66547.2
$ws.Cells.Item(108, 9).Value = 50000  # I108: 0 -> 50000
$ws.Cells.Item(108, 10).Value = 70684  # J108: 70347.2 -> 70684
$ws.Cells.Item(108, 11).Value = 50000  # K108: 0 -> 50000
$ws.Cells.Item(108, 12).Value = 70684  # L108: 70347.2 -> 70684
$ws.Cells.Item(108, 14).Value = -78364  # N108: -78027.2 -> -78364
$ws.Cells.Item(108, 13).Value = -46160  # M108: None -> -46160
$ws.Cells.Item(109, 8).Value = 65142.5  # H109: 70285 -> 65142.5
$ws.Cells.Item(109, 10).Value = 65142.5  # J109: 70285 -> 65142.5
$ws.Cells.Item(109, 12).Value = 65142.5  # L109: 70285 -> 65142.5
$ws.Cells.Item(109, 14).Value = -67222.5  # N109: -72365 -> -67222.5
$ws.Cells.Item(111, 8).Value = 78361.60000000001  # H111: 80702 -> 78361.60000000001
$ws.Cells.Item(111, 10).Value = 78361.60000000001  # J111: 80702 -> 78361.60000000001
$ws.Cells.Item(111, 12).Value = 78361.60000000001  # L111: 80702 -> 78361.60000000001
$ws.Cells.Item(111, 14).Value = -86541.60000000001  # N111: -88882 -> -86541.60000000001
$ws.Cells.Item(136, 8).Value = 516815.22  # H136: 563631.4 -> 516815.22
$ws.Cells.Item(136, 9).Value = 687667.5  # I136: 687660.5600000001 -> 687667.5
$ws.Cells.Item(136, 10).Value = 4258.3335  # J136: 5500 -> 4258.3335
$ws.Cells.Item(136, 11).Value = 2063002.5  # K136: 2062981.68 -> 2063002.5
$ws.Cells.Item(136, 12).Value = 12775.0005  # L136: 16500 -> 12775.0005
$ws.Cells.Item(136, 13).Value = -2060452.5  # M136: -2060431.68 -> -2060452.5
$ws.Cells.Item(136, 14).Value = -17875.0005  # N136: -21600 -> -17875.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 45458030  # H17: 3651.2856 -> 45458030
$ws.Cells.Item(17, 10).Value = 250019000  # J17: 25339.334 -> 250019000
$ws.Cells.Item(17, 12).Value = 750057000  # L17: 76018.00199999999 -> 750057000
$ws.Cells.Item(17, 14).Value = -750057338  # N17: -76356.00199999999 -> -750057338
$ws.Cells.Item(22, 8).Value = 399.66666  # H22: 288.33334 -> 399.66666
$ws.Cells.Item(22, 9).Value = 399.5  # I22: 463 -> 399.5
$ws.Cells.Item(22, 10).Value = 400  # J22: 201 -> 400
$ws.Cells.Item(22, 11).Value = 1198.5  # K22: 1389 -> 1198.5
$ws.Cells.Item(22, 12).Value = 1200  # L22: 603 -> 1200
$ws.Cells.Item(22, 13).Value = -1029.5  # M22: -1220 -> -1029.5
$ws.Cells.Item(22, 14).Value = -1538  # N22: -941 -> -1538
$ws.Cells.Item(27, 8).Value = 399.66666  # H27: 288.33334 -> 399.66666
$ws.Cells.Item(27, 9).Value = 399.5  # I27: 463 -> 399.5
$ws.Cells.Item(27, 10).Value = 400  # J27: 201 -> 400
$ws.Cells.Item(27, 11).Value = 1198.5  # K27: 1389 -> 1198.5
$ws.Cells.Item(27, 12).Value = 1200  # L27: 603 -> 1200
$ws.Cells.Item(27, 13).Value = -1096.5  # M27: -1287 -> -1096.5
$ws.Cells.Item(27, 14).Value = -1404  # N27: -807 -> -1404
$ws.Cells.Item(81, 8).Value = 5754.8  # H81: 100005704 -> 5754.8
$ws.Cells.Item(81, 9).Value = 3497.6667  # I81: 333336670 -> 3497.6667
$ws.Cells.Item(81, 11).Value = 10493.0001  # K81: 1000010010 -> 10493.0001
$ws.Cells.Item(81, 13).Value = -9370.000100000001  # M81: -1000008887 -> -9370.000100000001
$ws.Cells.Item(84, 8).Value = 5754.8  # H84: 100005704 -> 5754.8
$ws.Cells.Item(84, 9).Value = 3497.6667  # I84: 333336670 -> 3497.6667
$ws.Cells.Item(84, 11).Value = 31479.0003  # K84: 3000030030 -> 31479.0003
$ws.Cells.Item(84, 13).Value = -25863.0003  # M84: -3000024414 -> -25863.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3203.7083  # H113: 2899.875 -> 3203.7083
$ws.Cells.Item(113, 9).Value = 1720.6364  # I113: 1588.3077 -> 1720.6364
$ws.Cells.Item(113, 10).Value = 4458.615  # J113: 4449.909 -> 4458.615
$ws.Cells.Item(113, 11).Value = 1720.6364  # K113: 1588.3077 -> 1720.6364
$ws.Cells.Item(113, 12).Value = 4458.615  # L113: 4449.909 -> 4458.615
$ws.Cells.Item(113, 13).Value = 449.3635999999999  # M113: 581.6922999999999 -> 449.3635999999999
$ws.Cells.Item(113, 14).Value = -8798.615  # N113: -8789.909 -> -8798.615

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1445.7142  # H55: 885.9091 -> 1445.7142
$ws.Cells.Item(55, 9).Value = 674  # I55: 770 -> 674
$ws.Cells.Item(55, 10).Value = 2024.5  # J55: 982.5 -> 2024.5
$ws.Cells.Item(55, 11).Value = 674  # K55: 770 -> 674
$ws.Cells.Item(55, 12).Value = 2024.5  # L55: 982.5 -> 2024.5
$ws.Cells.Item(55, 13).Value = -501  # M55: -597 -> -501
$ws.Cells.Item(55, 14).Value = -2370.5  # N55: -1328.5 -> -2370.5
$ws.Cells.Item(68, 8).Value = 5339.737  # H68: 5442.5 -> 5339.737
$ws.Cells.Item(68, 10).Value = 5820.643  # J68: 5999.923 -> 5820.643
$ws.Cells.Item(68, 12).Value = 5820.643  # L68: 5999.923 -> 5820.643
$ws.Cells.Item(68, 14).Value = -7318.643  # N68: -7497.923 -> -7318.643
$ws.Cells.Item(71, 8).Value = 5339.737  # H71: 5442.5 -> 5339.737
$ws.Cells.Item(71, 10).Value = 5820.643  # J71: 5999.923 -> 5820.643
$ws.Cells.Item(71, 12).Value = 29103.215  # L71: 29999.615 -> 29103.215
$ws.Cells.Item(71, 14).Value = -36591.215  # N71: -37487.615 -> -36591.215
$ws.Cells.Item(100, 8).Value = 9651.5  # H100: 9023.277 -> 9651.5
$ws.Cells.Item(100, 9).Value = 2342.8572  # I100: 2566.6667 -> 2342.8572
$ws.Cells.Item(100, 10).Value = 15336  # J100: 12251.583 -> 15336
$ws.Cells.Item(100, 11).Value = 2342.8572  # K100: 2566.6667 -> 2342.8572
$ws.Cells.Item(100, 12).Value = 15336  # L100: 12251.583 -> 15336
$ws.Cells.Item(100, 13).Value = -1801.8572  # M100: -2025.6667 -> -1801.8572
$ws.Cells.Item(100, 14).Value = -16418  # N100: -13333.583 -> -16418
$ws.Cells.Item(122, 8).Value = 4625.846  # H122: 4714.5835 -> 4625.846
$ws.Cells.Item(122, 9).Value = 4027.7856  # I122: 4165 -> 4027.7856
$ws.Cells.Item(122, 10).Value = 5323.5835  # J122: 5364.091 -> 5323.5835
$ws.Cells.Item(122, 11).Value = 12083.3568  # K122: 12495 -> 12083.3568
$ws.Cells.Item(122, 12).Value = 15970.7505  # L122: 16092.273 -> 15970.7505
$ws.Cells.Item(122, 13).Value = -9633.356800000001  # M122: -10045 -> -9633.356800000001
$ws.Cells.Item(122, 14).Value = -20870.7505  # N122: -20992.273 -> -20870.7505
$ws.Cells.Item(132, 8).Value = 2778.6482  # H132: 2879.6924 -> 2778.6482
$ws.Cells.Item(132, 9).Value = 2143.8164  # I132: 2228.5957 -> 2143.8164
$ws.Cells.Item(132, 11).Value = 6431.449200000001  # K132: 6685.7871 -> 6431.449200000001
$ws.Cells.Item(132, 13).Value = -3901.449200000001  # M132: -4155.7871 -> -3901.449200000001
$ws.Cells.Item(136, 8).Value = 4024.4897  # H136: 4180.6 -> 4024.4897
$ws.Cells.Item(136, 9).Value = 3340.476  # I136: 3541.718 -> 3340.476
$ws.Cells.Item(136, 10).Value = 8128.5713  # J136: 8333.333000000001 -> 8128.5713
$ws.Cells.Item(136, 11).Value = 10021.428  # K136: 10625.154 -> 10021.428
$ws.Cells.Item(136, 12).Value = 24385.7139  # L136: 24999.999 -> 24385.7139
$ws.Cells.Item(136, 13).Value = -7471.428  # M136: -8075.153999999999 -> -7471.428
$ws.Cells.Item(136, 14).Value = -29485.7139  # N136: -30099.999 -> -29485.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 2149.4443  # H113: 2293.125 -> 2149.4443
$ws.Cells.Item(113, 10).Value = 2260  # J113: 2575 -> 2260
$ws.Cells.Item(113, 12).Value = 6780  # L113: 7725 -> 6780
$ws.Cells.Item(113, 14).Value = -11120  # N113: -12065 -> -11120
$ws.Cells.Item(136, 8).Value = 9903347  # H136: 9600731 -> 9903347
$ws.Cells.Item(136, 9).Value = 11179956  # I136: 10799620 -> 11179956
$ws.Cells.Item(136, 11).Value = 33539868  # K136: 32398860 -> 33539868
$ws.Cells.Item(136, 13).Value = -33537318  # M136: -32396310 -> -33537318
